## Progress.xlsx update: "Update Progress and add Idioms and Words"
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Progress"
$ws2 = $wb.Worksheets.Item(2)   # "TPO Results"

# ---------------------------------------------------------------------
# 1. Progress updates on already-existing rows
# ---------------------------------------------------------------------
$ws1.Range("D9").Value  = 234
$ws1.Range("D10").Value = 434
$ws1.Range("D16").Value = 5
$ws1.Range("C20").Value = 52
$ws1.Range("D20").Value = 42

# ---------------------------------------------------------------------
# 2. Make room for 5 new books: push the current row 38 (and the blank
#    spacer row 39 below it) down to rows 43/44, keeping their content
#    and formatting intact.
# ---------------------------------------------------------------------
$ws1.Rows("38:42").Insert()

# Grow Table1 (book list) so the new rows are part of the table/autofilter
$tbl = $ws1.ListObjects.Item(1)
$tbl.Resize($ws1.Range("A2:J43"))

# ---------------------------------------------------------------------
# 3. Fill in the 5 new rows. Cell values are entered in the same order
#    the original author typed them so new shared strings line up.
# ---------------------------------------------------------------------

# Row 38 - Vocabulary and Grammar for the TOEFL - Collins
$ws1.Range("A38").Value = 36
$ws1.Range("B38").Value = "Vocabulary and Grammar for the TOEFL - Collins"
$ws1.Range("J38").Value = "Review the copied pages "
$ws1.Range("C38").Value = 170
$ws1.Range("D38").Value = 108
$ws1.Range("F38").Value = "Y"
$ws1.Range("I38").Value = "Y"

# Row 40 - Word Power Made Easy (typed before row 39 by the author)
$ws1.Range("A40").Value = 38
$ws1.Range("B40").Value = "Word Power Made Easy"
$ws1.Range("C40").Value = 528
$ws1.Range("D40").Value = 70
$ws1.Range("F40").Value = "Y"
$ws1.Range("I40").Value = "Y"

# Row 39 - TED-Ed videos
$ws1.Range("A39").Value = 37
$ws1.Range("B39").Value = "TED-Ed videos"
$ws1.Range("C39").Value = 100
$ws1.Range("D39").Value = 0
$ws1.Range("G39").Value = "Y"
$ws1.Range("H39").Value = "Y"

# Row 41 - BBC Learning English Towards Advanced Course
$ws1.Range("A41").Value = 39
$ws1.Range("B41").Value = "BBC Learning English Towards Advanced Course"
$ws1.Range("C41").Value = 30
$ws1.Range("D41").Value = 2
$ws1.Range("F41").Value = "Y"
$ws1.Range("G41").Value = "Y"
$ws1.Range("H41").Value = "Y"
$ws1.Range("I41").Value = "Y"

# Row 42 - Grammar Booster
$ws1.Range("A42").Value = 40
$ws1.Range("B42").Value = "Grammar Booster"
$ws1.Range("C42").Value = 120
$ws1.Range("D42").Value = 10
$ws1.Range("H42").Value = "Y"
$ws1.Range("I42").Value = "Y"

# Row 43 already holds the previous row-38 data (TOEFL Grammar Tests Book)
# after the insert above; just make sure its calculated "Progress" cell
# keeps working by restoring the table formula explicitly.
$ws1.Range("E43").Formula = ' CONCATENATE(CEILING(Table1[[#This Row],[Current Part/Page]]/Table1[[#This Row],[Part/Page count]], 0.0001) * 100,"%")'

# ---------------------------------------------------------------------
# 4. TPO Results: fill in score for entry 41
# ---------------------------------------------------------------------
$ws2.Range("B42").Value = 23

# ---------------------------------------------------------------------
# 5. Selections, matching the saved view state
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D41").Select()

$ws2.Activate()
$ws2.Range("C42").Select()

$ws1.Activate()
